$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New summary rows 75-78 (text labels first, in the order new shared strings should appear)
$ws.Range("E75").Value = "median"
$ws.Range("E76").Value = "average"
$ws.Range("E77").Value = "max"
$ws.Range("E78").Value = "min"

# Update header for column F (shared string "ann" -> "ann  (ave quarters)")
$ws.Range("F1").Value = "ann  (ave quarters)"

# Column F: replace hard-coded averages with AVERAGE(B:E) formulas
$ws.Range("F2").Formula = "=AVERAGE(B2:E2)"
$ws.Range("F3:F66").Formula = "=AVERAGE(B3:E3)"
$ws.Range("F67:F72").Formula = "=AVERAGE(B67:E67)"

$ws.Range("F75").Formula = "=MEDIAN(F2:F72)"
$ws.Range("F76").Formula = "=AVERAGE(F2:F72)"
$ws.Range("F77").Formula = "=MAX(F2:F72)"
$ws.Range("F78").Formula = "=MIN(F2:F72)"

# Make sure all of column F (rows with stray "0%" style) use the consistent "0.00%" style
$ws.Range("F27").NumberFormat = "0.00%"
$ws.Range("F29").NumberFormat = "0.00%"
$ws.Range("F37").NumberFormat = "0.00%"
$ws.Range("F59").NumberFormat = "0.00%"
$ws.Range("F63").NumberFormat = "0.00%"
$ws.Range("F71").NumberFormat = "0.00%"
$ws.Range("F77").NumberFormat = "0.00%"
$ws.Range("F78").NumberFormat = "0.00%"

# Column width for F
$ws.Columns("F").ColumnWidth = 17.85546875

# Update sheet view selection
$ws.Range("F79").Select()

Write-Host "done"
